$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C13").Value = 0.000449900793650797
$ws.Range("D13").Value = 0.0633333333333333

$ws.Range("C14").Value = 0.00165555555555555
$ws.Range("D14").Value = 0.0853333333333333

$ws.Range("C15").Value = 0.208966397849462
